$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 590 (pushes old rows 590..631 down to 591..632),
# matching the diff's expansion of the used range from A1:D631 to A1:D632.
$ws.Rows(590).Insert()

# Write the new row's values. The date-like text "2026/01/08" must stay a
# plain text string (like the rest of column A) rather than being
# auto-converted to a serial date number by the smart-entry parser, so we
# round-trip it through a text formula cell and paste only the value back,
# which preserves the literal string without forcing a new number-format
# style onto the cell.
$ws.Range("ZZ1").Formula = "=""2026/01/08"""
$ws.Range("ZZ1").Copy()
$ws.Range("A590").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

$ws.Range("B590").Value = "木"
$ws.Range("C590").Value = 6
$ws.Range("D590").Value = 201
